$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-27 and 41-50 price/coin/link/volume updates per 2022-12-18 refresh

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '252.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.544'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05691'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.453'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8069'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.040'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1431'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03154'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02941'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09275'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001667'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.215'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04783'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005813'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16OneONE'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006459'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005056'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("B20").Value = 'UpBots'
$ws.Range("C20").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.007491'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19UpBotsUBXTBestin24h'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001051'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001501'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.986'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'GateToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.381'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '23GateTokenGT'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.090'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3320'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1277'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26ProBitTokenPROB'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006963'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003502'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1045'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009556'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005643'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7856'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.01703'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01011'
$ws.Range("D50").Style = "Normal"
